$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Set the new value for E9 (the additional control in the feed)
$ws.Range("E9").Value = 1

# Update the selection to E10, as recorded in the saved view state
$ws.Range("E10").Select()
